$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (updates workbook.xml sheet name)
$ws.Name = "Through 2022-02-25"

# Update the shared string / label for row 3 (February)
$ws.Range("A3").Value = "February (through 02-25)"

# Update February row (row 3) values
$ws.Range("C3").Value = 31
$ws.Range("E3").Value = 49
$ws.Range("F3").Value = 27
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 110
$ws.Range("I3").Value = 125

# Update Total row (row 4) values
$ws.Range("C4").Value = 82
$ws.Range("E4").Value = 135
$ws.Range("F4").Value = 76
$ws.Range("G4").Value = 138
$ws.Range("H4").Value = 327
$ws.Range("I4").Value = 284
